# Apply the change described by the diff:
# On the "About" worksheet (sheet1), add a new cell C1 containing the date
# 2021-04-21 (Excel serial 44307) formatted as a date (numFmtId 14).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

$cell = $ws.Range("C1")
$cell.Value = Get-Date -Year 2021 -Month 4 -Day 21 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$cell.NumberFormat = "mm-dd-yy"
